$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry text (index 61 in sharedStrings, appended)
$desc = "1. MiT module working understanding - completed`n2. MiT B0 kaggle training 2ep tried : low mIoU maybe due to no pretrained weights`n3. Alternate model definition tried - similar results`n4. LR_finder tried to get new LR - 5e-5`n5. MiT-B1 model trained for 2ep with LR_finder, similar results`n6.  Need to check Poly LR scheduler with 12ep training, else huggingface transformer implementation check"

# New row 52 data (Sno 51, date 12-Jul-2022, 02:00 -> 05:20)
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = 44754
$ws.Cells.Item(52, 3).Value = 0.08333333333333333
$ws.Cells.Item(52, 4).Value = 0.2222222222222222
$ws.Cells.Item(52, 5).Formula = "=D52-C52"
$ws.Cells.Item(52, 6).Value = "Code"
$ws.Cells.Item(52, 7).Value = $desc

# Match styling of the row above (number formats, alignment, wrap text)
$ws.Range("A51:G51").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)
$ws.Cells.Item(52, 7).Value = $desc
$excel.CutCopyMode = 0

# Row height for the wrapped, 6-line (8 visually wrapped) description
$ws.Rows.Item(52).RowHeight = 120

# Update total formula range and recalc
$ws.Cells.Item(57, 5).Formula = "=SUM(E2:E53)"

# Restore view state to match post-edit scroll/selection
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Range("G53").Select()
